$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append new data row (row 11) mirroring the existing "Room 101" / "none" rows,
# with the new measure_name "No more problems" and spc_chart_type "t".
$ws.Cells.Item(11, 1).Value = 101
$ws.Cells.Item(11, 2).Value = "No more problems"
$ws.Cells.Item(11, 3).Value = "Room 101"
$ws.Cells.Item(11, 4).Value = "t"
$ws.Cells.Item(11, 5).Value = "none"

# Slightly widen column A (ref column) to fit the new content.
$ws.Columns.Item(1).ColumnWidth = 3.2

# Move/leave the active selection where the author last left it.
$ws.Range("B14").Select()
